$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Control Lines")
$ws2 = $wb.Worksheets.Item("Test Program")

# --- "Test Program" sheet: insert a new (binary opcode) column before the
# existing data, shifting old B->C, C->D, D->E.
$ws2.Range("B1").EntireColumn.Insert()

# Row 5 ("LD A, [addr]"): fix the operand-format cell (old value "00" is
# replaced by the new "0[addr]" value, now living in D5).
$ws2.Range("D5").Value = "0[addr]"

# New row 6: "JP [addr]" instruction encoding.
$ws2.Range("A6").Value = "JP [addr]"
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "001011"
$ws2.Range("C6").Value = "2c"
$ws2.Range("D6").Value = "0[addr]"
$ws2.Range("E6").Value = "[addr]"

# --- Sheet view / selection state ---
# "Test Program" is no longer the active tab; its own cell-selection moves
# on to D7.
$ws2.Range("D7").Select()

# "Control Lines" becomes the active tab, with X6 selected.
$ws1.Activate()
$ws1.Range("X6").Select()
